$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Module 3 practical session moves from Thursday 20.02 (Aud J) to Monday 24.02 (Aud C).
# Row 9, column C (Monday) currently holds "24.02: No lecture" -> becomes the practical session.
$ws.Range("C9").Value = "24.02: <strong>Practical session</strong>  in Aud C."

# Row 8, column D (Thursday) currently holds the practical session text -> becomes "No lecture."
$ws.Range("D8").Value = "20.02: No lecture."

# Move the active selection to D9, matching the saved workbook state.
$ws.Range("D9").Select()
